# 2023 Attendance workbook update:
#  - Ohio North SQL Saturday (Cleveland) 2023, row 10: fill in Registered/Attended
#  - SQL Saturday Denver 2023, row 11: fill in Attended count
#  - New event added as row 16: SQL Saturday Los Angeles 2023

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: fill in Registered (D) / Attended (E) and compute the No-show rate (F)
$ws.Range("D10").Value = 100
$ws.Range("E10").Value = 70
$ws.Range("F10").Formula = "=IF(D10=0,0,+(D10-E10)/D10)"

# Row 11: fill in Attended (E) count
$ws.Range("E11").Value = 120

# Row 16: brand new event entry
$ws.Range("A16").Value = "SQL Saturday Los Angeles 2023"
$ws.Range("B16").Value = 1049
$ws.Range("C16").Value = "6/10/2023"
$ws.Range("D16").Value = 257
$ws.Range("E16").Value = 135
$ws.Range("F16").Formula = "=IF(D16=0,0,+(D16-E16)/D16)"

# Leave the selection where it was left in the saved file
$ws.Range("F10").Select() | Out-Null
